# "Update countries & provincias Spain" — refresh the COVID country table on
# sheet "Pais" with the newer snapshot figures and re-rank a handful of
# countries whose totals just crossed a neighbour's.
#
# Note on shared strings: a handful of rows keep the SAME row number but now
# show a DIFFERENT country, because in the source workbook the country list
# got re-sorted by "Casos totales" and two rows traded places. We reproduce
# that by writing the new country name (and its own updated stats) into the
# row that used to hold the other country — Excel resolves `.Value = "..."`
# against the existing shared-string table on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 16:05"

# --- Pure stat refreshes (country/ranking unchanged) -------------------
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1387499
$ws.Range("C4").Value = 1665
$ws.Range("E4").Value = 1043337
$ws.Range("G4").Value = 142
$ws.Range("H4").Value = 81937

# Row 11 (Brasil)
$ws.Range("B11").Value = 170021
$ws.Range("C11").Value = 878
$ws.Range("E11").Value = 90936
$ws.Range("G11").Value = 76
$ws.Range("H11").Value = 11701

# Row 48 (Singapur)
$ws.Range("B48").Value = 10243
$ws.Range("C48").Value = 67
$ws.Range("D48").Value = 3600
$ws.Range("E48").Value = 6423
$ws.Range("F48").Value = 23
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 220

# Row 75 (Republica de Macedonia)
$ws.Range("B75").Value = 2519
$ws.Range("C75").Value = 33
$ws.Range("D75").Value = 2010
$ws.Range("E75").Value = 499

# Row 81 (Hong Kong)
$ws.Range("B81").Value = 1995
$ws.Range("C81").Value = 109
$ws.Range("D81").Value = 742
$ws.Range("E81").Value = 1234

# Row 82 (Kirguistan)
$ws.Range("D82").Value = 1776
$ws.Range("E82").Value = 15

# Row 86 (Maldivas)
$ws.Range("B86").Value = 1674
$ws.Range("C86").Value = 10
$ws.Range("D86").Value = 1205
$ws.Range("E86").Value = 377
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 92

# --- Sri Lanka overtakes Libano (rows 105/106 swap) ---------------------
$ws.Range("A105").Value = "Sri Lanka"
$ws.Range("B105").Value = 872
$ws.Range("C105").Value = 9
$ws.Range("D105").Value = 366
$ws.Range("E105").Value = 497
$ws.Range("F105").Value = 1
$ws.Range("H105").Value = 9

$ws.Range("A106").Value = "Libano"
$ws.Range("B106").Value = 870
$ws.Range("C106").Value = 11
$ws.Range("D106").Value = 234
$ws.Range("E106").Value = 610
$ws.Range("F106").Value = 3
$ws.Range("H106").Value = 26

# --- Zambia jumps up to row 125, shifting Taiwan..Cabo Verde down one ---
$ws.Range("A125").Value = "Zambia"
$ws.Range("B125").Value = 441
$ws.Range("C125").Value = 174
$ws.Range("D125").Value = 117
$ws.Range("E125").Value = 317
$ws.Range("F125").Value = 1

$ws.Range("A126").Value = "Taiwan"
$ws.Range("B126").Value = 440
$ws.Range("D126").Value = 372
$ws.Range("E126").Value = 61
$ws.Range("H126").Value = 7

$ws.Range("A127").Value = "Guinea Ecuatorial"
$ws.Range("B127").Value = 439
$ws.Range("D127").Value = 13
$ws.Range("E127").Value = 422
$ws.Range("F127").Value = 0
$ws.Range("H127").Value = 4

$ws.Range("A128").Value = "Reunion"
$ws.Range("B128").Value = 436
$ws.Range("D128").Value = 354
$ws.Range("E128").Value = 82
$ws.Range("F128").Value = 5
$ws.Range("H128").Value = 0

$ws.Range("A129").Value = "Venezuela"
$ws.Range("B129").Value = 422
$ws.Range("D129").Value = 205
$ws.Range("E129").Value = 207
$ws.Range("F129").Value = 2
$ws.Range("H129").Value = 10

$ws.Range("A130").Value = "Estado de Palestina"
$ws.Range("B130").Value = 375
$ws.Range("D130").Value = 316
$ws.Range("E130").Value = 57
$ws.Range("H130").Value = 2

$ws.Range("A131").Value = "Sierra Leona"
$ws.Range("B131").Value = 338
$ws.Range("D131").Value = 72
$ws.Range("E131").Value = 247
$ws.Range("H131").Value = 19

$ws.Range("A132").Value = "Congo"
$ws.Range("B132").Value = 333
$ws.Range("D132").Value = 53
$ws.Range("E132").Value = 269
$ws.Range("H132").Value = 11

$ws.Range("A133").Value = "Mauricio"
$ws.Range("B133").Value = 332
$ws.Range("D133").Value = 322
$ws.Range("E133").Value = 0
$ws.Range("F133").Value = 0
$ws.Range("H133").Value = 10

$ws.Range("A134").Value = "Isla de Man"
$ws.Range("B134").Value = 330
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 271
$ws.Range("E134").Value = 36
$ws.Range("F134").Value = 21
$ws.Range("H134").Value = 23

$ws.Range("A135").Value = "Benin"
$ws.Range("B135").Value = 327
$ws.Range("C135").Value = 8
$ws.Range("D135").Value = 76
$ws.Range("E135").Value = 249
$ws.Range("F135").Value = 0
$ws.Range("H135").Value = 2

$ws.Range("A136").Value = "Montenegro"
$ws.Range("B136").Value = 324
$ws.Range("D136").Value = 294
$ws.Range("E136").Value = 21
$ws.Range("F136").Value = 2
$ws.Range("H136").Value = 9

$ws.Range("A137").Value = "Republica del Chad"
$ws.Range("B137").Value = 322
$ws.Range("D137").Value = 53
$ws.Range("E137").Value = 238
$ws.Range("F137").Value = 0
$ws.Range("H137").Value = 31

$ws.Range("A138").Value = "Vietnam"
$ws.Range("B138").Value = 288
$ws.Range("D138").Value = 252
$ws.Range("E138").Value = 36
$ws.Range("F138").Value = 2

$ws.Range("A139").Value = "Ruanda"
$ws.Range("B139").Value = 285
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 150
$ws.Range("E139").Value = 135
$ws.Range("H139").Value = 0

$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("C140").Value = 7
$ws.Range("D140").Value = 58
$ws.Range("E140").Value = 207
$ws.Range("F140").Value = 0
$ws.Range("H140").Value = 2

# --- San Bartolome overtakes Sahara Occidental (rows 215/216 swap) ------
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"
